# "fixed for red line"
# Insert a new stop row ("STATION; GLENBURY IN") right before the existing
# "STATION; OVERBROOK OUT" row (row 10) in the Green Line schedule, shifting
# all subsequent rows down by one. The new row reuses the same column
# layout/style as the rest of the stop rows (B/C/D, one per train column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 10 (and everything below it) down by one row.
$ws.Rows.Item(10).Insert()

# Fill the newly-opened row with the new station stop for all three trains.
$ws.Range("B10:D10").Value = "STATION; GLENBURY IN"

# Match the author's last on-screen selection after making the edit.
[void]$ws.Range("E9").Select()
